$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The last paragraph holds the breakout repo URL, wrongly auto-capitalised
# (e.g. "Https://github.Com/..."). Fix the casing to a proper lowercase URL.
# Editing the run directly (rather than the whole paragraph/text range)
# keeps it as a single run so its existing character formatting (colour,
# cap="none", ...) is preserved instead of being duplicated/reset.
$para = $tr.Paragraphs($tr.Paragraphs().Count)
$run1 = $para.Runs(1)
$run1.Text = "https://github.com/seanyoung247/breakout"
